$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder: 10/12/2015 -> 10/13/2015
#    The date footer placeholder (type dt = 16) appears once on the slide
#    master and once on every slide layout; update every occurrence.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "10/12/2015") {
                $tr.Text = "10/13/2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 7, "TextBox 2": reword "laser peak" -> "laser position" and split
#    the run so "laser " / "position " stand on their own, then grow the
#    auto-fit textbox to match the new (wider) text.
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$shape = $slide7.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# "peak" -> "position"
$peakIdx = $tr.Text.IndexOf("peak")
if ($peakIdx -ge 0) {
    $tr.Characters($peakIdx + 1, 4).Text = "position"
}

# Break "laser " into its own run.
$laserIdx = $tr.Text.IndexOf("laser ")
if ($laserIdx -ge 0) {
    $tr.Characters($laserIdx + 1, 6).Text = "laser "
}

# Break "position " into its own run (shifts the trailing space off "in image:").
$positionIdx = $tr.Text.IndexOf("position ")
if ($positionIdx -ge 0) {
    $tr.Characters($positionIdx + 1, 9).Text = "position "
}

# The text box auto-fits its width to the text (wrap="none" + spAutoFit);
# grow it to the new rendered width (5898025 EMU == 464.41141732283467 pt).
# Shape.Width is a single-precision COM property, so nudge past the f32
# rounding boundary to land exactly on the target EMU value.
$shape.Width = (5898025 / 12700.0) + 0.00001
